$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item(1).Name = "GNG_TO-16509961211945608"
$wb.Worksheets.Item(2).Name = "NB_TO-16509961225305672"
$wb.Worksheets.Item(3).Name = "RS_TO-16509961225305672"
$wb.Worksheets.Item(4).Name = "TOL_TO-16509961225945983"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16509961226585956"

# Sheet 1 (GNG) - update B2:B5
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16509961211545763.csv"
$ws1.Range("B3").Value = "GNG_stims-16509961211786077.csv"
$ws1.Range("B4").Value = "go_stims-16509961211786077.csv"
$ws1.Range("B5").Value = "GNG_stims-16509961211945608.csv"

# Sheet 2 (NB) - update B2:B10
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-1650996122434561.csv"
$ws2.Range("B3").Value = "OB-16509961220185628.csv"
$ws2.Range("B4").Value = "ZB-match_3-16509961212905655.csv"
$ws2.Range("B5").Value = "TB-16509961221145825.csv"
$ws2.Range("B6").Value = "OB-1650996122002567.csv"
$ws2.Range("B7").Value = "TB-16509961224985604.csv"
$ws2.Range("B8").Value = "ZB-match_4-16509961217466078.csv"
$ws2.Range("B9").Value = "OB-16509961218505616.csv"
$ws2.Range("B10").Value = "ZB-match_4-16509961213545692.csv"

# Sheet 4 (TOL) - update B2:B7
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16509961225626192.csv"
$ws4.Range("B3").Value = "ZM_stims-16509961225385664.csv"
$ws4.Range("B4").Value = "MM_stims-16509961225785637.csv"
$ws4.Range("B5").Value = "ZM_stims-16509961225626192.csv"
$ws4.Range("B6").Value = "MM_stims-16509961225945983.csv"
$ws4.Range("B7").Value = "ZM_stims-16509961225785637.csv"

# Sheet 5 (vSAT) - update B2:B5
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16509961225945983.csv"
$ws5.Range("B3").Value = "SAT_stims-16509961226105957.csv"
$ws5.Range("B4").Value = "vSAT_stims-16509961226265619.csv"
$ws5.Range("B5").Value = "vSAT_stims-16509961226426032.csv"
